$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.607.80'
$ws.Range("E2").Value = '  -2.31%  '

$ws.Range("D3").Value = '1.806.47'
$ws.Range("E3").Value = '  -1.75%  '

$ws.Range("E4").Value = '  +0.44%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '230.02'
$ws.Range("E5").Value = '  -0.53%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.611'
$ws.Range("E6").Value = '  +0.04%  '

$ws.Range("E7").Value = '  +0.46%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '39.51'
$ws.Range("E8").Value = '  -9.91%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.323'
$ws.Range("E9").Value = '  +4.50%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0680'
$ws.Range("E10").Value = '  -3.45%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0991'
$ws.Range("E11").Value = '  -1.73%  '

$ws.Range("D12").Value = '2.070.57'
$ws.Range("E12").Value = '  -1.58%  '

$ws.Range("E13").Value = '  -0.25%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.663'
$ws.Range("E14").Value = '  -1.52%  '

$ws.Range("D15").Value = '1.803.02'
$ws.Range("E15").Value = '  -1.92%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '4.58'
$ws.Range("E16").Value = '  -2.49%  '

$ws.Range("D17").Value = '34.665.77'
$ws.Range("E17").Value = '  -2.09%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '69.35'
$ws.Range("E18").Value = '  -1.01%  '

$ws.Range("D19").Value = '0.0₃0783'
$ws.Range("E19").Value = '  -2.24%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '239.72'
$ws.Range("E20").Value = '  -1.93%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.87'
$ws.Range("E21").Value = '  -1.96%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.66'
$ws.Range("E22").Value = '  -0.76%  '

$ws.Range("E23").Value = '  +0.42%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.23'
$ws.Range("E24").Value = '  +1.69%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '171.79'
$ws.Range("E25").Value = '  +0.71%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.76'
$ws.Range("E26").Value = '  -1.87%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.18'
$ws.Range("E27").Value = '  -3.11%  '

$ws.Range("E28").Value = '  +0.14%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.49'
$ws.Range("E29").Value = '  -3.12%  '

$ws.Range("E30").Value = '  +0.28%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.05'
$ws.Range("E31").Value = '  +2.98%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0545'
$ws.Range("E32").Value = '  -1.73%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.92'
$ws.Range("E33").Value = '  -4.20%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.27'
$ws.Range("E34").Value = '  +15.16%  '

$ws.Range("E35").Value = '  -3.33%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.696'
$ws.Range("E36").Value = '  +1.39%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '91.50'
$ws.Range("E37").Value = '  -4.58%  '

$ws.Range("E38").Value = '  +4.49%  '

$ws.Range("D39").Value = '1.324.62'
$ws.Range("E39").Value = '  -1.59%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0191'
$ws.Range("E40").Value = '  -1.77%  '

$ws.Range("E41").Value = '  +0.35%  '

$ws.Range("B42").Value = 'InjectiveProtocol'
$ws.Range("C42").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '14.45'
$ws.Range("E42").Value = '  -7.60%  '

$ws.Range("B43").Value = 'ARBITRUM'
$ws.Range("C43").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.958'
$ws.Range("E43").Value = '  -5.28%  '

$ws.Range("B44").Value = 'MXToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.74'
$ws.Range("E44").Value = '  -2.42%  '

$ws.Range("B45").Value = 'RenderToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.20'
$ws.Range("E45").Value = '  -9.99%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '6.22'
$ws.Range("E46").Value = '  -0.68%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0512'
$ws.Range("E47").Value = '  -1.46%  '

$ws.Range("D48").Value = '1.995.37'
$ws.Range("E48").Value = '  -0.51%  '

$ws.Range("E49").Value = '  +0.36%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0664'
$ws.Range("E50").Value = '  +5.53%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '97.96'
$ws.Range("E51").Value = '  -4.91%  '
